$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 98.62294467674243
$ws.Range("R2").Value = 887.6065020906819
$ws.Range("S2").Value = 0.005018734788036125
$ws.Range("T2").Value = 0.005018734788036125
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 258.9360906292527
$ws.Range("R3").Value = 2330.424815663274
$ws.Range("S3").Value = 0.01317676703102503
$ws.Range("T3").Value = 0.01317676703102503
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 18.58532251067512
$ws.Range("R4").Value = 167.267902596076
$ws.Range("S4").Value = 0.0009457718478891906
$ws.Range("T4").Value = 0.0009457718478891904
$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 6.996142233409778
$ws.Range("R5").Value = 62.96528010068801
$ws.Range("S5").Value = 0.0003560204222653126
$ws.Range("T5").Value = 0.0003560204222653126
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 4661.668733572595
$ws.Range("R6").Value = 41955.01860215335
$ws.Range("S6").Value = 0.2372234891197517
$ws.Range("T6").Value = 0.2372234891197517
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.6228339974986082
$ws.Range("T7").Value = 0.6228339974986082
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 878.4833705316089
$ws.Range("R8").Value = 7906.35033478448
$ws.Range("S8").Value = 0.04470435421340577
$ws.Range("T8").Value = 0.04470435421340577
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 330.6907699015822
$ws.Range("R9").Value = 2976.21692911424
$ws.Range("S9").Value = 0.01682822670147775
$ws.Range("T9").Value = 0.01682822670147775
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 258.408016945683
$ws.Range("R10").Value = 2325.672152511147
$ws.Range("S10").Value = 0.01314989436183974
$ws.Range("T10").Value = 0.01314989436183974
$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 678.454308117531
$ws.Range("R11").Value = 6106.08877305778
$ws.Range("S11").Value = 0.03452525423371795
$ws.Range("T11").Value = 0.03452525423371795
$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 48.696541661994
$ws.Range("R12").Value = 438.2688749579461
$ws.Range("S12").Value = 0.002478074737041733
$ws.Range("T12").Value = 0.002478074737041732
$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 18.331020704472
$ws.Range("R13").Value = 164.979186340248
$ws.Range("S13").Value = 0.0009328309108117666
$ws.Range("T13").Value = 0.0009328309108117666
$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 39.58912301041255
$ws.Range("R14").Value = 356.3021070937129
$ws.Range("S14").Value = 0.002014615458212474
$ws.Range("T14").Value = 0.002014615458212474
$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 103.9418644145823
$ws.Range("R15").Value = 935.476779731241
$ws.Range("S15").Value = 0.005289404535431775
$ws.Range("T15").Value = 0.005289404535431775
$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 7.46050142261489
$ws.Range("R16").Value = 67.144512803534
$ws.Range("S16").Value = 0.0003796507815559081
$ws.Range("T16").Value = 0.000379650781555908
$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 2.808384361110222
$ws.Range("R17").Value = 25.275459249992
$ws.Range("S17").Value = 0.0001429133589295909
$ws.Range("T17").Value = 0.0001429133589295909
